$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 972.9
$ws.Range("I15").Value = 972.9
$ws.Range("K15").Value = 2918.7
$ws.Range("M15").Value = -2749.7

$ws.Range("H76").Value = 6105
$ws.Range("I76").Value = 5121.6665
$ws.Range("J76").Value = 6432.778
$ws.Range("K76").Value = 5121.6665
$ws.Range("L76").Value = 6432.778
$ws.Range("M76").Value = -4806.6665
$ws.Range("N76").Value = -7062.778

$ws.Range("H79").Value = 6105
$ws.Range("I79").Value = 5121.6665
$ws.Range("J79").Value = 6432.778
$ws.Range("K79").Value = 5121.6665
$ws.Range("L79").Value = 6432.778
$ws.Range("M79").Value = -4029.6665
$ws.Range("N79").Value = -8616.778

$ws.Range("H88").Value = 3005549.2
$ws.Range("I88").Value = 161
$ws.Range("J88").Value = 3270730.5
$ws.Range("K88").Value = 161
$ws.Range("L88").Value = 3270730.5
$ws.Range("M88").Value = 245
$ws.Range("N88").Value = -3271542.5

$ws.Range("H91").Value = 3005549.2
$ws.Range("I91").Value = 161
$ws.Range("J91").Value = 3270730.5
$ws.Range("K91").Value = 161
$ws.Range("L91").Value = 3270730.5
$ws.Range("M91").Value = 1243
$ws.Range("N91").Value = -3273538.5

$ws.Range("H97").Value = 2614.1428
$ws.Range("J97").Value = 2614.1428
$ws.Range("L97").Value = 7842.428400000001
$ws.Range("N97").Value = -8834.428400000001

$ws.Range("H99").Value = 1091.3846
$ws.Range("I99").Value = 263.75
$ws.Range("J99").Value = 1459.2222
$ws.Range("K99").Value = 791.25
$ws.Range("L99").Value = 4377.6666
$ws.Range("M99").Value = 706.75
$ws.Range("N99").Value = -7373.6666

$ws.Range("H104").Value = 245.66667
$ws.Range("I104").Value = 245.66667
$ws.Range("K104").Value = 737.00001
$ws.Range("M104").Value = 1009.99999

$ws.Range("H106").Value = 1881.5
$ws.Range("I106").Value = 1553.3334
$ws.Range("K106").Value = 1553.3334
$ws.Range("M106").Value = -922.3334

$ws.Range("H123").Value = 120000
$ws.Range("J123").Value = 120000
$ws.Range("L123").Value = 120000
$ws.Range("N123").Value = -129800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6853599.5
$ws.Range("I45").Value = 9591894
$ws.Range("K45").Value = 9591894
$ws.Range("M45").Value = -9591517

$ws.Range("H46").Value = 3510.2
$ws.Range("J46").Value = 3137.75
$ws.Range("L46").Value = 3137.75
$ws.Range("N46").Value = -3775.75

$ws.Range("H63").Value = 4613.48
$ws.Range("I63").Value = 2625.6924
$ws.Range("K63").Value = 2625.6924
$ws.Range("M63").Value = -1939.6924

$ws.Range("H66").Value = 4613.48
$ws.Range("I66").Value = 2625.6924
$ws.Range("K66").Value = 13128.462
$ws.Range("M66").Value = -9696.462

$ws.Range("H74").Value = 56690.53
$ws.Range("I74").Value = 3219.25
$ws.Range("K74").Value = 3219.25
$ws.Range("M74").Value = -2345.25

$ws.Range("H77").Value = 56690.53
$ws.Range("I77").Value = 3219.25
$ws.Range("K77").Value = 16096.25
$ws.Range("M77").Value = -11728.25

$ws.Range("H122").Value = 1101718.5
$ws.Range("I122").Value = 3666
$ws.Range("K122").Value = 10998
$ws.Range("M122").Value = -8548

$ws.Range("H132").Value = 3936.1875
$ws.Range("I132").Value = 2436.75
$ws.Range("K132").Value = 7310.25
$ws.Range("M132").Value = -4780.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 59000
$ws.Range("J53").Value = 59000
$ws.Range("L53").Value = 59000
$ws.Range("N53").Value = -60148

$ws.Range("H58").Value = 14910
$ws.Range("J58").Value = 14910
$ws.Range("L58").Value = 14910
$ws.Range("N58").Value = -15498

$ws.Range("H59").Value = 122923
$ws.Range("J59").Value = 122923
$ws.Range("L59").Value = 122923
$ws.Range("N59").Value = -124617

$ws.Range("H134").Value = 4383
$ws.Range("I134").Value = 2145.2632
$ws.Range("J134").Value = 10456.857
$ws.Range("K134").Value = 6435.7896
$ws.Range("L134").Value = 31370.571
$ws.Range("M134").Value = -3900.7896
$ws.Range("N134").Value = -36440.571

$ws.Range("H141").Value = 103270
$ws.Range("J141").Value = 103270
$ws.Range("L141").Value = 103270
$ws.Range("N141").Value = -113630

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1845.0769
$ws.Range("I16").Value = 1781
$ws.Range("J16").Value = 2197.5
$ws.Range("K16").Value = 1781
$ws.Range("L16").Value = 2197.5
$ws.Range("M16").Value = -1494
$ws.Range("N16").Value = -2771.5

$ws.Range("H58").Value = 5505.4688
$ws.Range("I58").Value = 5987.273
$ws.Range("J58").Value = 4445.5
$ws.Range("K58").Value = 5987.273
$ws.Range("L58").Value = 4445.5
$ws.Range("M58").Value = -5784.273
$ws.Range("N58").Value = -4851.5

$ws.Range("H62").Value = 3778.6
$ws.Range("I62").Value = 2200
$ws.Range("K62").Value = 2200
$ws.Range("M62").Value = -1576

$ws.Range("H65").Value = 3778.6
$ws.Range("I65").Value = 2200
$ws.Range("K65").Value = 11000
$ws.Range("M65").Value = -7880

$ws.Range("H86").Value = 13423.143
$ws.Range("I86").Value = 9808
$ws.Range("K86").Value = 9808
$ws.Range("M86").Value = -8685

$ws.Range("H89").Value = 13423.143
$ws.Range("I89").Value = 9808
$ws.Range("K89").Value = 49040
$ws.Range("M89").Value = -43424

$ws.Range("H107").Value = 2149.48
$ws.Range("I107").Value = 1760.1
$ws.Range("J107").Value = 3707
$ws.Range("K107").Value = 1760.1
$ws.Range("L107").Value = 3707
$ws.Range("M107").Value = 159.9000000000001
$ws.Range("N107").Value = -7547

$ws.Range("H113").Value = 1845.0769
$ws.Range("I113").Value = 1781
$ws.Range("J113").Value = 2197.5
$ws.Range("K113").Value = 1781
$ws.Range("L113").Value = 2197.5
$ws.Range("M113").Value = 389
$ws.Range("N113").Value = -6537.5

$ws.Range("H134").Value = 2835.932
$ws.Range("I134").Value = 1936.8857
$ws.Range("J134").Value = 6332.222
$ws.Range("K134").Value = 5810.6571
$ws.Range("L134").Value = 18996.666
$ws.Range("M134").Value = -3275.6571
$ws.Range("N134").Value = -24066.666

$ws.Range("H136").Value = 5505.4688
$ws.Range("I136").Value = 5987.273
$ws.Range("J136").Value = 4445.5
$ws.Range("K136").Value = 17961.819
$ws.Range("L136").Value = 13336.5
$ws.Range("M136").Value = -15411.819
$ws.Range("N136").Value = -18436.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 995
$ws.Range("I59").Value = 995
$ws.Range("K59").Value = 2985
$ws.Range("M59").Value = -2445

$ws.Range("H98").Value = 1734.8572
$ws.Range("I98").Value = 1062.3334
$ws.Range("J98").Value = 1918.2727
$ws.Range("K98").Value = 3187.0002
$ws.Range("L98").Value = 5754.8181
$ws.Range("M98").Value = -1689.0002
$ws.Range("N98").Value = -8750.8181

$ws.Range("H103").Value = 433
$ws.Range("I103").Value = 149.5
$ws.Range("K103").Value = 448.5
$ws.Range("M103").Value = 430.5

$ws.Range("H131").Value = 18945924
$ws.Range("J131").Value = 19615180
$ws.Range("L131").Value = 58845540
$ws.Range("N131").Value = -58855620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9755

$ws.Range("H29").Value = 9659.666999999999
$ws.Range("I29").Value = 8237.5
$ws.Range("J29").Value = 12504
$ws.Range("K29").Value = 8237.5
$ws.Range("L29").Value = 12504
$ws.Range("M29").Value = -7947.5
$ws.Range("N29").Value = -13084

$ws.Range("H51").Value = 61400
$ws.Range("J51").Value = 61400
$ws.Range("L51").Value = 61400
$ws.Range("N51").Value = -62418

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H126").Value = 4027363.2
$ws.Range("I126").Value = 3248782.5
$ws.Range("J126").Value = 4632926
$ws.Range("K126").Value = 9746347.5
$ws.Range("L126").Value = 13898778
$ws.Range("M126").Value = -9743877.5
$ws.Range("N126").Value = -13903718

$ws.Range("H132").Value = 3422.2593
$ws.Range("I132").Value = 2993.7778
$ws.Range("K132").Value = 8981.3334
$ws.Range("M132").Value = -6451.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 9000
$ws.Range("I23").Value = 9000
$ws.Range("K23").Value = 9000
$ws.Range("M23").Value = -8770

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H61").Value = 10103452
$ws.Range("I61").Value = 11112197
$ws.Range("K61").Value = 11112197
$ws.Range("M61").Value = -11111995

$ws.Range("H113").Value = 10103452
$ws.Range("I113").Value = 11112197
$ws.Range("K113").Value = 11112197
$ws.Range("M113").Value = -11110027

$ws.Range("H122").Value = 7670.364
$ws.Range("J122").Value = 9347.5
$ws.Range("L122").Value = 28042.5
$ws.Range("N122").Value = -32942.5

$ws.Range("H138").Value = 82498.336
$ws.Range("J138").Value = 82498.336
$ws.Range("L138").Value = 82498.336
$ws.Range("N138").Value = -92778.336
